$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '45.133.63'
$ws.Range('E2').Value = '  -3.38%  '
$ws.Range('D3').Value = '2.380.51'
$ws.Range('E3').Value = '  +4.66%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.17%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '293.06'
$ws.Range('E5').Value = '  -3.17%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '94.47'
$ws.Range('E6').Value = '  -6.16%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.556'
$ws.Range('E7').Value = '  -1.27%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '1.00'
$ws.Range('E8').Value = '  +0.08%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.498'
$ws.Range('E9').Value = '  -3.98%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '34.01'
$ws.Range('E10').Value = '  -5.43%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0774'
$ws.Range('E11').Value = '  -1.73%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '6.93'
$ws.Range('E12').Value = '  -4.02%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.103'
$ws.Range('E13').Value = '  +0.71%  '
$ws.Range('D14').Value = '2.744.47'
$ws.Range('D15').Value = '2.377.05'
$ws.Range('E15').Value = '  +4.54%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '13.91'
$ws.Range('E16').Value = '  +1.58%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.820'
$ws.Range('E17').Value = '  +2.11%  '
$ws.Range('D18').Value = '45.121.81'
$ws.Range('E18').Value = '  -3.34%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.36'
$ws.Range('E19').Value = '  -6.04%  '
$ws.Range('D20').Value = '0.0₃0926'
$ws.Range('E20').Value = '  -0.51%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.06'
$ws.Range('E21').Value = '  +2.01%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '66.18'
$ws.Range('E22').Value = '  +1.11%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '237.77'
$ws.Range('E23').Value = '  -4.07%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.74'
$ws.Range('E24').Value = '  -4.04%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.999'
$ws.Range('E25').Value = '  -0.09%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.87'
$ws.Range('E26').Value = '  -0.78%  '
$ws.Range('B27').Value = 'Toncoin'
$ws.Range('C27').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.21'
$ws.Range('E27').Value = '  -1.61%  '
$ws.Range('B28').Value = 'InjectiveProtocol'
$ws.Range('C28').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '37.57'
$ws.Range('E28').Value = '  -12.03%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.49'
$ws.Range('E29').Value = '  -2.60%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '3.80'
$ws.Range('E30').Value = '  +15.26%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '20.93'
$ws.Range('E31').Value = '  +5.10%  '
$ws.Range('E32').Value = '  -2.24%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '147.00'
$ws.Range('E33').Value = '  -0.23%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.38'
$ws.Range('E34').Value = '  -2.61%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0755'
$ws.Range('E35').Value = '  -3.04%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.95'
$ws.Range('E36').Value = '  +11.76%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.111'
$ws.Range('E37').Value = '  -3.66%  '
$ws.Range('E38').Value = '  -2.02%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '14.70'
$ws.Range('E39').Value = '  -8.03%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.72'
$ws.Range('E40').Value = '  -5.14%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0293'
$ws.Range('E41').Value = '  -1.57%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.16'
$ws.Range('E42').Value = '  -3.36%  '
$ws.Range('D43').Value = '1.930.66'
$ws.Range('E43').Value = '  +6.06%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.00'
$ws.Range('E44').Value = '  +0.18%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '89.28'
$ws.Range('E45').Value = '  -1.33%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.71'
$ws.Range('E46').Value = '  -13.55%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '8.46'
$ws.Range('E47').Value = '  +8.08%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '14.92'
$ws.Range('E48').Value = '  +15.41%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '99.09'
$ws.Range('E49').Value = '  +4.57%  '
$ws.Range('D50').Value = '2.616.06'
$ws.Range('E50').Value = '  +4.68%  '
$ws.Range('E51').Value = '  -5.02%  '
